# Apply crypto price/volume updates from the latest scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.985.11'
$ws.Range('E2').Value = '  -5.36%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.824.09'
$ws.Range('E3').Value = '  -4.24%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.70%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.09'
$ws.Range('E5').Value = '  -3.25%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  -0.49%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4631'
$ws.Range('E7').Value = '  -2.96%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3845'
$ws.Range('E8').Value = '  -4.12%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07865'
$ws.Range('E9').Value = '  -2.29%  '

# Row 10
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9585'
$ws.Range('E10').Value = '  -3.39%  '

# Row 11
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.86'
$ws.Range('E11').Value = '  -5.88%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.798.20'
$ws.Range('E12').Value = '  -5.46%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.643'
$ws.Range('E13').Value = '  -5.08%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.860'
$ws.Range('E14').Value = '  -3.34%  '

# Row 15
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06845'
$ws.Range('E15').Value = '  +0.37%  '

# Row 16
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.68%  '

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '86.26'
$ws.Range('E17').Value = '  -3.37%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009926'
$ws.Range('E18').Value = '  -2.96%  '

# Row 19
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.62'
$ws.Range('E19').Value = '  -4.25%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9995'
$ws.Range('E20').Value = '  -0.55%  '

# Row 21
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.007.97'
$ws.Range('E21').Value = '  -5.29%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.301'
$ws.Range('E22').Value = '  -3.88%  '

# Row 23
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.95'
$ws.Range('E23').Value = '  -6.16%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.098'
$ws.Range('E24').Value = '  -2.60%  '

# Row 25
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.022.42'
$ws.Range('E25').Value = '  -4.36%  '

# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.37'
$ws.Range('E26').Value = '  -3.15%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.13'
$ws.Range('E27').Value = '  -2.82%  '

# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.730'
$ws.Range('E28').Value = '  -11.73%  '

# Row 29
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.967'
$ws.Range('E29').Value = '  -4.49%  '

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.44'
$ws.Range('E30').Value = '  -2.56%  '

# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9336'
$ws.Range('E31').Value = '  -6.27%  '

# Row 32
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09218'
$ws.Range('E32').Value = '  -3.66%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.279'
$ws.Range('E33').Value = '  -3.60%  '

# Row 34
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.315'
$ws.Range('E34').Value = '  -5.47%  '

# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.339'
$ws.Range('E35').Value = '  -5.62%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05925'
$ws.Range('E36').Value = '  -7.94%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02137'
$ws.Range('E37').Value = '  -4.99%  '

# Row 38
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.142'
$ws.Range('E38').Value = '  -4.45%  '

# Row 39
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9987'
$ws.Range('E39').Value = '  -0.79%  '

# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.557'
$ws.Range('E40').Value = '  -2.15%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5572'
$ws.Range('E41').Value = '  -4.48%  '

# Row 42
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.893'
$ws.Range('E42').Value = '  -6.31%  '

# Row 43
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1761'
$ws.Range('E43').Value = '  -3.37%  '

# Row 44
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.230'
$ws.Range('E44').Value = '  -3.48%  '

# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.195'
$ws.Range('E45').Value = '  -9.89%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.64'
$ws.Range('E46').Value = '  -4.60%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5236'
$ws.Range('E47').Value = '  -4.51%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06993'
$ws.Range('E48').Value = '  -5.53%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.820'
$ws.Range('E49').Value = '  -7.25%  '

# Row 50
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '111.60'
$ws.Range('E50').Value = '  -4.03%  '

# Row 51
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9988'
$ws.Range('E51').Value = '  -0.69%  '
